$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-10 (Efna4-Epha4 pair), columns E..T
# Ligand-expressing cell count for "ECs" (sending cluster) rose from 1 to 2,
# and Receptor-expressing cell count for "ECs" (target cluster) rose from 2 to 3,
# which cascades through all derived specificity / weight columns.

$values = @{
  2  = @(2, 0.6666666666666666, 0.2988413333333333, 0.896524, 0.3632971504731247, 0.3632971504731246, 3, 1, 8.081040666666667, 24.243122, 0.4661250698616886, 0.4661250698616886, 2.414948967547556, 21.734540707928, 0.1693419096448376, 0.1693419096448376)
  3  = @(2, 0.6666666666666666, 0.2988413333333333, 0.896524, 0.3632971504731247, 0.3632971504731246, 3, 1, 7.000300666666668, 21.000902, 0.4037865631294714, 0.4037865631294715, 2.09197918496089, 18.827812664648, 0.1466945077842734, 0.1466945077842734)
  4  = @(2, 0.6666666666666666, 0.2988413333333333, 0.896524, 0.3632971504731247, 0.3632971504731246, 3, 1, 2.255294666666666, 6.765884, 0.1300883670088399, 0.1300883670088399, 0.6739752652462222, 6.065777387215999, 0.04726073304401358, 0.04726073304401357)
  5  = @(2, 0.6666666666666666, 0.3851916666666667, 1.155575, 0.4682720202225272, 0.4682720202225272, 3, 1, 8.081040666666667, 24.243122, 0.4661250698616886, 0.4661250698616886, 3.112749522794444, 28.01474570515, 0.2182733281404996, 0.2182733281404996)
  6  = @(2, 0.6666666666666666, 0.3851916666666667, 1.155575, 0.4682720202225272, 0.4682720202225272, 3, 1, 7.000300666666668, 21.000902, 0.4037865631294714, 0.4037865631294715, 2.696457480961111, 24.26811732865, 0.1890819496553486, 0.1890819496553486)
  7  = @(2, 0.6666666666666666, 0.3851916666666667, 1.155575, 0.4682720202225272, 0.4682720202225272, 3, 1, 2.255294666666666, 6.765884, 0.1300883670088399, 0.1300883670088399, 0.8687207114777776, 7.8184864033, 0.06091674242667903, 0.06091674242667903)
  8  = @(2, 0.6666666666666666, 0.138548, 0.415644, 0.1684308293043481, 0.1684308293043481, 3, 1, 8.081040666666667, 24.243122, 0.4661250698616886, 0.4661250698616886, 1.119612022285333, 10.076508200568, 0.07850983207635143, 0.07850983207635143)
  9  = @(2, 0.6666666666666666, 0.138548, 0.415644, 0.1684308293043481, 0.1684308293043481, 3, 1, 7.000300666666668, 21.000902, 0.4037865631294714, 0.4037865631294715, 0.9698776567653336, 8.728898910888002, 0.06801010568984939, 0.0680101056898494)
  10 = @(2, 0.6666666666666666, 0.138548, 0.415644, 0.1684308293043481, 0.1684308293043481, 3, 1, 2.255294666666666, 6.765884, 0.1300883670088399, 0.1300883670088399, 0.3124665654773333, 2.812199089296, 0.02191089153814731, 0.02191089153814731)
}

$columns = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($row in $values.Keys) {
  $rowValues = $values[$row]
  for ($i = 0; $i -lt $columns.Length; $i++) {
    $col = $columns[$i]
    $ws.Range("$col$row").Value = $rowValues[$i]
  }
}
